$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '46.098.65'
$ws.Range("E2").Value = '  +3.44%  '
$ws.Range("D3").Value = '2.447.14'
$ws.Range("E3").Value = '  +0.86%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '319.87'
$ws.Range("E5").Value = '  +2.14%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '104.66'
$ws.Range("E6").Value = '  +2.90%  '
$ws.Range("E7").Value = '  +0.98%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("E9").Value = '  +4.39%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.89'
$ws.Range("E10").Value = '  +1.82%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0806'
$ws.Range("E11").Value = '  +0.40%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.122'
$ws.Range("E12").Value = '  -3.04%  '
$ws.Range("E13").Value = '  -3.53%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.06'
$ws.Range("E14").Value = '  +1.68%  '
$ws.Range("D15").Value = '2.827.53'
$ws.Range("E15").Value = '  +0.80%  '
$ws.Range("D16").Value = '2.430.62'
$ws.Range("E16").Value = '  +0.45%  '
$ws.Range("E17").Value = '  +0.66%  '
$ws.Range("D18").Value = '45.917.07'
$ws.Range("E18").Value = '  +3.25%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.52'
$ws.Range("E19").Value = '  +0.12%  '
$ws.Range("E20").Value = '  +0.27%  '
$ws.Range("E21").Value = '  -0.32%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '71.20'
$ws.Range("E22").Value = '  +3.25%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.38'
$ws.Range("E23").Value = '  +4.85%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '246.68'
$ws.Range("E24").Value = '  +2.34%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.50'
$ws.Range("E25").Value = '  +1.46%  '
$ws.Range("E26").Value = '  -0.02%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.86'
$ws.Range("E27").Value = '  +2.72%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.29'
$ws.Range("E28").Value = '  +1.04%  '
$ws.Range("E29").Value = '  +0.35%  '
$ws.Range("E30").Value = '  +1.23%  '
$ws.Range("E31").Value = '  +1.55%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.127'
$ws.Range("E32").Value = '  +2.37%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.10'
$ws.Range("E33").Value = '  +2.77%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.34'
$ws.Range("E34").Value = '  +2.92%  '
$ws.Range("E35").Value = '  -0.05%  '
$ws.Range("E36").Value = '  -1.17%  '
$ws.Range("E37").Value = '  -0.22%  '
$ws.Range("E38").Value = '  -0.13%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.95'
$ws.Range("E39").Value = '  +1.70%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '126.42'
$ws.Range("E40").Value = '  +1.80%  '
$ws.Range("E41").Value = '  +2.45%  '
$ws.Range("E42").Value = '  +1.29%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '20.87'
$ws.Range("E43").Value = '  -2.39%  '
$ws.Range("E44").Value = '  +1.20%  '
$ws.Range("D45").Value = '1.960.04'
$ws.Range("E45").Value = '  +0.53%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.96'
$ws.Range("E46").Value = '  +0.57%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.08'
$ws.Range("E47").Value = '  -4.60%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.86'
$ws.Range("E48").Value = '  +11.95%  '
$ws.Range("E49").Value = '  -4.52%  '
$ws.Range("E50").Value = '  +7.53%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '77.61'
$ws.Range("E51").Value = '  +4.99%  '
